$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "B2" = 1161.56531337
    "C2" = 1290.366832755
    "D2" = 1128.38323257
    "B3" = 1318.93042116
    "C3" = 1476.90389847
    "D3" = 1304.22379527
    "B4" = 1838.26286061
    "C4" = 2037.872878195
    "D4" = 1774.65915594
    "B5" = 1740.92691168
    "C5" = 1934.319344265
    "D5" = 1693.155997665
    "B6" = 1614.08635509
    "C6" = 1793.372514395
    "D6" = 1569.74176161
    "B7" = 1496.48059374
    "C7" = 1662.69101383
    "D7" = 1455.3419562
    "B8" = 1387.73169942
    "C8" = 1541.87615736
    "D8" = 1349.60810787
    "B9" = 1286.56100265
    "C9" = 1429.48568917
    "D9" = 1251.235521855
    "B10" = 1192.73770539
    "C10" = 1325.236894905
    "D10" = 1159.980749235
    "B11" = 1105.9239444
    "C11" = 1228.773574875
    "D11" = 1075.549397175
    "B12" = 1025.20965129
    "C12" = 1139.07572471
    "D12" = 997.020076635
    "B13" = 985.78139979
    "C13" = 1095.045481145
    "D13" = 957.4192232550001
    "B14" = 977.8035728999999
    "C14" = 1086.19635971
    "D14" = 949.7120138550001
    "B15" = 970.11390432
    "C15" = 1077.63008969
    "D15" = 942.204296385
    "B16" = 962.4954292799999
    "C16" = 1069.2015873
    "D16" = 934.86753063
    "B17" = 954.5737789499999
    "C17" = 1060.379130195
    "D17" = 927.112512645
    "B18" = 947.0004030299999
    "C18" = 1052.01493936
    "D18" = 919.89609426
    "B19" = 939.46930077
    "C19" = 1043.573234935
    "D19" = 912.3845083650001
    "B20" = 932.0807126399999
    "C20" = 1035.41694872
    "D20" = 905.328139815
    "B21" = 924.3897874199999
    "C21" = 1026.857358065
    "D21" = 897.8312118600001
    "B22" = 917.0888007599999
    "C22" = 1018.75616689
    "D22" = 890.76821724
    "B23" = 909.79125369
    "C23" = 1010.62405906
    "D23" = 883.59581475
    "B24" = 902.50135932
    "C24" = 1002.54305998
    "D24" = 876.560725755
    "B25" = 895.20702579
    "C25" = 994.468441505
    "D25" = 869.5521058950001
    "B26" = 888.1406522999999
    "C26" = 986.5680782650001
    "D26" = 862.57529028
    "B27" = 881.0896142399999
    "C27" = 978.7516457950001
    "D27" = 855.7481765250001
    "B28" = 873.8985614699999
    "C28" = 970.801453545
    "D28" = 848.85102621
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
